$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-14 Monday" "2024-10-15 Tuesday"

Replace-Text "276÷7=" "484÷5="
Replace-Text "854÷8=" "314÷4="
Replace-Text "811÷4=" "382÷3="
Replace-Text "346÷7=" "254÷2="
Replace-Text "963÷5=" "295÷5="

Replace-Text "931÷8=" "577÷9="
Replace-Text "143÷8=" "980÷8="
Replace-Text "715÷5=" "998÷4="
Replace-Text "234÷8=" "627÷7="
Replace-Text "563÷7=" "927÷3="

Replace-Text "188÷8=" "168÷4="
Replace-Text "278÷4=" "343÷7="
Replace-Text "724÷5=" "696÷7="
Replace-Text "757÷4=" "497÷6="
Replace-Text "609÷5=" "123÷2="

Replace-Text "219÷2=" "829÷5="
Replace-Text "781÷6=" "946÷9="
Replace-Text "259÷3=" "852÷8="
Replace-Text "321÷6=" "779÷6="
Replace-Text "360÷7=" "508÷7="

Replace-Text "690÷3=" "151÷3="
Replace-Text "812÷9=" "829÷7="
Replace-Text "333÷2=" "579÷9="
Replace-Text "223÷2=" "680÷8="
Replace-Text "841÷2=" "546÷8="
